# Update the Summary sheet values
$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = "Noor Al Suwaidi"
$summary.Range("B4").Value = 6506.58
$summary.Range("B6").Value = 4829
$summary.Range("B7").Value = 23603
$summary.Range("B8").Value = -18774
$summary.Range("B9").Value = 0.2

# Assets sheet: remove the "Vehicles / Premium Car" row (row 2), shifting
# "Liquid Assets / Savings Account" and the TOTAL row up by one.
$assets = $wb.Worksheets.Item("Assets")
$assets.Rows.Item(2).Delete()
$assets.Range("C2").Value = 4829
$assets.Range("C3").Value = 4829

# Liabilities sheet: remove the "Auto Loans / Vehicle Loan 1" row (row 2),
# shifting "Credit Cards / Credit Card Balance" and the TOTAL row up by one.
$liabilities = $wb.Worksheets.Item("Liabilities")
$liabilities.Rows.Item(2).Delete()
$liabilities.Range("C2").Value = 23603
$liabilities.Range("D2").Value = 1180
$liabilities.Range("C3").Value = 23603
